$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet "Sheet" -> "Sheet1" ---
$ws.Name = "Sheet1"

# --- Header row values (row 1, columns A..X) ---
$headers = @(
    "Date",
    "Model Name",
    "Exact Precision (Micro Avg)",
    "Exact Recall (Micro Avg)",
    "Exact F1 Score (Micro Avg)",
    "Exact Precision (Macro Avg)",
    "Exact Recall (Macro Avg)",
    "Exact F1 Score (Macro Avg)",
    "Exact Precision (Weighted Avg)",
    "Exact Recall (Weighted Avg)",
    "Exact F1 Score (Weighted Avg)",
    "Partial Precision",
    "Partial Recall",
    "Partial F1 Score",
    "Partial TP",
    "Partial FP",
    "Partial FN",
    "Support",
    "Accuracy",
    "Result Link",
    "Stats Link",
    "No of GPU Used",
    "Power Consumption",
    "Unnamed: 23"
)

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Style the header row: bold font, thin box border, centered/top-aligned ---
$headerRange = $ws.Range("A1:X1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- Data row 2 ---
$row2 = @(
    "09/11/2025",
    "Qwen2.5-14B-Instruct",
    0.5085714285714286,
    0.2996632996632997,
    0.3771186440677966,
    0.2480829108953851,
    0.1352454290298497,
    0.1686706048410837,
    0.5612725352323296,
    0.2996632996632997,
    0.3819343482832793,
    0.5942857142857143,
    0.3513513513513514,
    0.4416135881104034,
    104,
    71,
    192,
    297,
    0.9542240862792091,
    "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-14B-Instruct_3_shot.txt",
    "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-14B-Instruct_3_shot.txt",
    "4 MLGPU",
    "0.044 kWh",
    609
)

for ($c = 1; $c -lt $row2.Length; $c++) {
    $ws.Cells.Item(2, $c + 1).Value = $row2[$c]
}

# --- Data row 3 ---
$row3 = @(
    "09/12/2025",
    "Qwen2.5-14B-Instruct",
    0.4204081632653061,
    0.3468013468013468,
    0.3800738007380073,
    0.5774682005374145,
    0.3393417531160006,
    0.4015189946505205,
    0.5718210822257008,
    0.3468013468013468,
    0.4093152891508639,
    0.5102040816326531,
    0.4222972972972973,
    0.4621072088724585,
    125,
    120,
    171,
    297,
    0.9440383463151588,
    "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-14B-Instruct_3_shot.txt",
    "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-14B-Instruct_3_shot.txt",
    "4 MLGPU",
    "0.024 kWh"
)

for ($c = 1; $c -lt $row3.Length; $c++) {
    $ws.Cells.Item(3, $c + 1).Value = $row3[$c]
}

# --- Dates in column A must stay literal text ("MM/DD/YYYY"), not auto-converted
#     to Excel date serials. Force text format before writing, then strip the
#     formatting back off so the cell keeps default (no) style like its siblings. ---
$dateCellA2 = $ws.Cells.Item(2, 1)
$dateCellA2.NumberFormat = "@"
$dateCellA2.Value = "09/11/2025"
$dateCellA2.ClearFormats()

$dateCellA3 = $ws.Cells.Item(3, 1)
$dateCellA3.NumberFormat = "@"
$dateCellA3.Value = "09/12/2025"
$dateCellA3.ClearFormats()
